$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 46, shifting the existing rows 46..153 down to 47..154.
$ws.Rows("46:46").Insert()

# Populate the newly inserted row 46 with the new price-report record.
$ws.Cells.Item(46, 1).Value = 10
$ws.Cells.Item(46, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(46, 3).Value = 'La Araucanía'
$ws.Cells.Item(46, 4).Value = 45281
$ws.Cells.Item(46, 5).Value = 9
$ws.Cells.Item(46, 6).Value = 'Fruta'
$ws.Cells.Item(46, 7).Value = 100108
$ws.Cells.Item(46, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(46, 9).Value = 100108004
$ws.Cells.Item(46, 10).Value = 'Papaya'
$ws.Cells.Item(46, 11).Value = 'Cultivar IV Región'
$ws.Cells.Item(46, 12).Value = 'Primera'
$ws.Cells.Item(46, 13).Value = 100
$ws.Cells.Item(46, 14).Value = 24000
$ws.Cells.Item(46, 15).Value = 24000
$ws.Cells.Item(46, 16).Value = 24000
$ws.Cells.Item(46, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(46, 18).Value = 'Provincia del Elquí'
$ws.Cells.Item(46, 19).Value = 2400
$ws.Cells.Item(46, 20).Value = 10
